$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price for existing rows (350 TL -> 300 TL, keeping trailing space)
$ws.Range("B2").Value = "300 TL "
$ws.Range("B3").Value = "300 TL "

# Add new row 4: Baggy Siyah
$ws.Range("A4").Value = "Baggy Siyah"
$ws.Range("B4").Value = "300 TL"
$ws.Range("C4").Value = "Jeans"
$ws.Range("D4").Value = "BAG6.jpg"
$ws.Range("E4").Value = "100% Pamuk"
$ws.Range("F4").Value = "Var"

# Update selection to F4
$ws.Range("F4").Select()
